$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arrow = [char]0x2190

$texts = @{
    2  = "$arrow That's me not caring."
    3  = "I pretend you are muted."
    4  = "I know you are doing Emails in parallel."
    5  = "I couldn't care less."
    6  = "Your hair looks funny."
    7  = "What's the purpose of this meeting again?"
    8  = "My silence isn't approval. I was just not listening."
    9  = "My video is not frozen. I just try not to move."
    10 = "I'm just trying to stick my cursor in your ear."
    11 = "We ignore your agenda."
    12 = "That will escalate quickly."
    13 = "I want to see the world burn."
    14 = 'You already lost me at "who sent the invitation for this?"'
    15 = "You switched off your camera so I moved a funny photo over your video thumbnail."
    16 = "Look there! A squirrel!"
    17 = "Here we are now, entertain us!"
    18 = 'All I hear is "mi mi mi".'
    19 = "Told you so."
    20 = "That will never scale."
    21 = "Who hired you?"
}

foreach ($row in $texts.Keys) {
    $ws.Range("B$row").Value = $texts[$row]
}

$ws.Range("B2").Select() | Out-Null
